$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Stephen Curry"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Golden State Warriors"

# Row 5
$ws.Range("A5").Value = "Austin Reaves"
$ws.Range("C5").Value = "Los Angeles Lakers"

# Row 12
$ws.Range("A12").Value = "Jalen Duren"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "Detroit Pistons"

# Row 13
$ws.Range("A13").Value = "Franz Wagner"
$ws.Range("B13").Value = "SF,PF"
$ws.Range("C13").Value = "Orlando Magic"

# Row 14
$ws.Range("A14").Value = "Darius Garland"
$ws.Range("B14").Value = "PG"
$ws.Range("C14").Value = "Cleveland Cavaliers"

# Row 15
$ws.Range("A15").Value = "Jaden McDaniels"
$ws.Range("C15").Value = "Minnesota Timberwolves"

# Row 16
$ws.Range("A16").Value = "Trey Murphy III"
$ws.Range("B16").Value = "SF,PF"
$ws.Range("C16").Value = "New Orleans Pelicans"
